{"js": "// Replace the date line and each \"a\u00f7b=c, r\" answer cell with its updated value.\n// Every \"before\" string in this worksheet is unique, so a single in-order\n// search+replace pass (first match each time) reproduces the diff exactly,\n// even though one new value (\"112\u00f72=56, 0\") coincides with an earlier original value.\nconst replacements = [\n  [\"2025-09-12 Friday\", \"2025-09-13 Saturday\"],\n  [\"112\u00f72=56, 0\", \"113\u00f77=16, 1\"],\n  [\"649\u00f78=81, 1\", \"562\u00f78=70, 2\"],\n  [\"747\u00f76=124, 3\", \"726\u00f78=90, 6\"],\n  [\"629\u00f79=69, 8\", \"663\u00f74=165, 3\"],\n  [\"886\u00f74=221, 2\", \"755\u00f79=83, 8\"],\n  [\"217\u00f79=24, 1\", \"910\u00f75=182, 0\"],\n  [\"510\u00f77=72, 6\", \"885\u00f74=221, 1\"],\n  [\"785\u00f72=392, 1\", \"962\u00f72=481, 0\"],\n  [\"673\u00f73=224, 1\", \"492\u00f74=123, 0\"],\n  [\"769\u00f72=384, 1\", \"183\u00f73=61, 0\"],\n  [\"529\u00f74=132, 1\", \"112\u00f72=56, 0\"],\n  [\"531\u00f73=177, 0\", \"615\u00f79=68, 3\"],\n  [\"165\u00f79=18, 3\", \"542\u00f72=271, 0\"],\n  [\"440\u00f76=73, 2\", \"617\u00f73=205, 2\"],\n  [\"400\u00f79=44, 4\", \"968\u00f78=121, 0\"],\n  [\"743\u00f78=92, 7\", \"333\u00f75=66, 3\"],\n  [\"793\u00f76=132, 1\", \"124\u00f76=20, 4\"],\n  [\"311\u00f72=155, 1\", \"655\u00f73=218, 1\"],\n  [\"621\u00f73=207, 0\", \"603\u00f74=150, 3\"],\n  [\"540\u00f79=60, 0\", \"158\u00f76=26, 2\"],\n  [\"453\u00f76=75, 3\", \"111\u00f76=18, 3\"],\n  [\"393\u00f78=49, 1\", \"978\u00f76=163, 0\"],\n  [\"477\u00f76=79, 3\", \"124\u00f74=31, 0\"],\n  [\"317\u00f75=63, 2\", \"580\u00f72=290, 0\"],\n  [\"227\u00f75=45, 2\", \"634\u00f77=90, 4\"],\n];\n\nfor (const [before, after] of replacements) {\n  const results = context.document.body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(after, Word.InsertLocation.replace);\n    await context.sync();\n  }\n}", "ps1": "# Update the date line and each \"a\u00f7b=c, r\" answer cell with its new value.\n# Every \"before\" string in this document is unique, so running one\n# Find/Replace (replace first match, in this document order) per pair\n# reproduces the diff exactly - even though one new value (\"112\u00f72=56, 0\")\n# happens to equal an original value used earlier in the document; by the\n# time we get to that later pair the earlier occurrence has already been\n# rewritten, so no double-replacement happens.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '2025-09-12 Friday'\n$find.Replacement.Text = '2025-09-13 Saturday'\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '112\u00f72=56, 0'\n$find.Replacement.Text = '113\u00f77=16, 1'\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '649\u00f78=81, 1'\n$find.Replacement.Text = '562\u00f78=70, 2'\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '747\u00f76=124, 3'\n$find.Replacement.Text = '726\u00f78=90, 6'\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '629\u00f79=69, 8'\n$find.Replacement.Text = '663\u00f74=165, 3'\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '886\u00f74=221, 2'\n$find.Replacement.Text = '755\u00f79=83, 8'\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '217\u00f79=24, 1'\n$find.Replacement.Text = '910\u00f75=182, 0'\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '510\u00f77=72, 6'\n$find.Replacement.Text = '885\u00f74=221, 1'\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '785\u00f72=392, 1'\n$find.Replacement.Text = '962\u00f72=481, 0'\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '673\u00f73=224, 1'\n$find.Replacement.Text = '492\u00f74=123, 0'\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '769\u00f72=384, 1'\n$find.Replacement.Text = '183\u00f73=61, 0'\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '529\u00f74=132, 1'\n$find.Replacement.Text = '112\u00f72=56, 0'\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '531\u00f73=177, 0'\n$find.Replacement.Text = '615\u00f79=68, 3'\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '165\u00f79=18, 3'\n$find.Replacement.Text = '542\u00f72=271, 0'\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '440\u00f76=73, 2'\n$find.Replacement.Text = '617\u00f73=205, 2'\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '400\u00f79=44, 4'\n$find.Replacement.Text = '968\u00f78=121, 0'\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '743\u00f78=92, 7'\n$find.Replacement.Text = '333\u00f75=66, 3'\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '793\u00f76=132, 1'\n$find.Replacement.Text = '124\u00f76=20, 4'\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '311\u00f72=155, 1'\n$find.Replacement.Text = '655\u00f73=218, 1'\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '621\u00f73=207, 0'\n$find.Replacement.Text = '603\u00f74=150, 3'\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '540\u00f79=60, 0'\n$find.Replacement.Text = '158\u00f76=26, 2'\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '453\u00f76=75, 3'\n$find.Replacement.Text = '111\u00f76=18, 3'\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '393\u00f78=49, 1'\n$find.Replacement.Text = '978\u00f76=163, 0'\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '477\u00f76=79, 3'\n$find.Replacement.Text = '124\u00f74=31, 0'\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '317\u00f75=63, 2'\n$find.Replacement.Text = '580\u00f72=290, 0'\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '227\u00f75=45, 2'\n$find.Replacement.Text = '634\u00f77=90, 4'\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n"}
